# 运维事项确认表.xlsx - update strategy so files and settings for Explorer
#
# Sheet "策略更新" (the first/active sheet) tracks strategy updates. This
# edit moves the tracked session from day-session ("日盘") to night-session
# ("夜盘") for both 郑州 (Zhengzhou) and 大连 (Dalian) rows, bumps the date
# from 2017-02-15 (serial 42781) to 2017-02-16 (serial 42782), rewrites the
# "内容" (content) notes for both rows with the new lot-size / so updates,
# and tweaks column widths / row heights so the taller two-line content
# fits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("策略更新")
$ws.Activate()

# --- Row 3: 郑州 (Zhengzhou) ---------------------------------------------

# C3: 日盘 -> 夜盘
$ws.Range("C3").Value = "夜盘"

# D3: date 2017-02-15 -> 2017-02-16 (Excel serial 42781 -> 42782)
$ws.Range("D3").Value = 42782

# E3: rewritten content, 3 runs: plain prefix, SimSun (宋体) label, Tahoma body
$e3Text = "(ok) 1. " + "更新郑州夜盘" + "so`n2. night`n(ok) zzsr fl34 4`n(ok) zzzc fl34 1`n(ok) zzcf fl34 5`n(ok) zzta fl34 10`n(ok) zzsr fl36 5`n(ok) zzcf fl36 5`n(ok) zzrm fl36 5`n(ok) zzta fl36 10`n(ok) zzcf fw10 4`n(ok) zzrm fw10 3`n(ok) zzsr fd10 5`n(ok) zzta fd10 5`n"
$e3 = $ws.Range("E3")
$e3.Value = $e3Text
$e3.Characters(9, 6).Font.Name = "宋体"
$e3.Characters(15, 218).Font.Name = "Tahoma"

# --- Row 4: 大连 (Dalian) -------------------------------------------------

# C4: 日盘 -> 夜盘
$ws.Range("C4").Value = "夜盘"

# D4: date 2017-02-15 -> 2017-02-16 (Excel serial 42781 -> 42782)
$ws.Range("D4").Value = 42782

# E4: rewritten content, 5 runs alternating plain/SimSun/Tahoma/SimSun/Tahoma
$e4Text = "(ok) 1. " + "更新大连夜盘" + "so`n2. " + "手数`n" + "(ok) dlm fl34 2`n(ok) dly fl34 3`n(ok) dla fl34 3`n(ok) dli fl34 10`n(ok) dla fl36 4`n(ok) dla fw10 3`n(ok) dlm fd10 5`n"
$e4 = $ws.Range("E4")
$e4.Value = $e4Text
$e4.Characters(9, 6).Font.Name = "宋体"
$e4.Characters(15, 6).Font.Name = "Tahoma"
$e4.Characters(21, 3).Font.Name = "宋体"
$e4.Characters(24, 113).Font.Name = "Tahoma"

# --- Column widths ---------------------------------------------------------
# Stored (OOXML) widths go 3 -> 5.25, 6.375 -> 8.125, 8.375 -> 10.375.
# ColumnWidth (COM, character units) is offset from the stored width by the
# sheet's standard 5px/MDW padding (5/7 for this workbook's Tahoma-11
# default font) and is itself pixel-quantized, so we feed it the nearest
# reachable character width for each target.
$ws.Columns.Item(1).ColumnWidth = 32 / 7   # -> stored width 37/7 = 5.2857 (target 5.25)
$ws.Columns.Item(2).ColumnWidth = 52 / 7   # -> stored width 57/7 = 8.1429 (target 8.125)
$ws.Columns.Item(3).ColumnWidth = 68 / 7   # -> stored width 73/7 = 10.4286 (target 10.375)

# --- Row heights -------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 26.25
$ws.Rows.Item(3).RowHeight = 52.5
$ws.Rows.Item(4).RowHeight = 52.5

# --- View: scroll the frozen-less sheet so column B is left-most visible --
$ws.Range("F4").Select()
try {
    $excel.ActiveWindow.ScrollColumn = 2
} catch {
}
